$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 2525287.5
$ws.Range("I31").Value = 2525287.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 7575862.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -7575632.5
$ws.Range("N31").ClearContents()

$ws.Range("H98").Value = 3730.5527
$ws.Range("I98").Value = 2298.8125
$ws.Range("J98").Value = 11366.5
$ws.Range("K98").Value = 2298.8125
$ws.Range("L98").Value = 11366.5
$ws.Range("M98").Value = -800.8125
$ws.Range("N98").Value = -14362.5

$ws.Range("H113").Value = 2605.4167
$ws.Range("I113").Value = 2008.125
$ws.Range("J113").Value = 3800
$ws.Range("K113").Value = 2008.125
$ws.Range("L113").Value = 3800
$ws.Range("M113").Value = 1245.875
$ws.Range("N113").Value = -10308

$ws.Range("H121").Value = 1351.2858
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1351.2858
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 4053.8574
$ws.Range("N121").Value = -7547.857400000001
$ws.Range("M121").ClearContents()

$ws.Range("H122").Value = 3730.5527
$ws.Range("I122").Value = 2298.8125
$ws.Range("J122").Value = 11366.5
$ws.Range("K122").Value = 6896.4375
$ws.Range("L122").Value = 34099.5
$ws.Range("M122").Value = -4446.4375
$ws.Range("N122").Value = -38999.5

$ws.Range("H138").Value = 2069
$ws.Range("I138").Value = 1346.7059
$ws.Range("J138").Value = 2216.9397
$ws.Range("K138").Value = 4040.1177
$ws.Range("L138").Value = 6650.8191
$ws.Range("M138").Value = 1099.8823
$ws.Range("N138").Value = -16930.8191

$ws.Range("H140").Value = 94667.375
$ws.Range("I140").Value = 37500
$ws.Range("J140").Value = 102834.14
$ws.Range("K140").Value = 37500
$ws.Range("L140").Value = 102834.14
$ws.Range("M140").Value = -32320
$ws.Range("N140").Value = -113194.14

$ws.Range("H141").Value = 7383.2666
$ws.Range("I141").Value = 4303.5
$ws.Range("J141").Value = 50500
$ws.Range("K141").Value = 12910.5
$ws.Range("L141").Value = 151500
$ws.Range("M141").Value = -7730.5
$ws.Range("N141").Value = -161860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30813.309
$ws.Range("I32").Value = 30045.84
$ws.Range("J32").Value = 50000
$ws.Range("K32").Value = 30045.84
$ws.Range("L32").Value = 50000
$ws.Range("M32").Value = -29758.84
$ws.Range("N32").Value = -50574

$ws.Range("H61").Value = 1452.3
$ws.Range("I61").Value = 1369.2222
$ws.Range("J61").Value = 2200
$ws.Range("K61").Value = 1369.2222
$ws.Range("L61").Value = 2200
$ws.Range("M61").Value = -1157.2222

$ws.Range("H74").Value = 992.0769
$ws.Range("I74").Value = 786.26666
$ws.Range("J74").Value = 1272.7273
$ws.Range("K74").Value = 786.26666
$ws.Range("L74").Value = 1272.7273
$ws.Range("M74").Value = 87.73334
$ws.Range("N74").Value = -3020.7273

$ws.Range("H77").Value = 992.0769
$ws.Range("I77").Value = 786.26666
$ws.Range("J77").Value = 1272.7273
$ws.Range("K77").Value = 3931.3333
$ws.Range("L77").Value = 6363.636500000001
$ws.Range("M77").Value = 436.6666999999998
$ws.Range("N77").Value = -15099.6365

$ws.Range("H136").Value = 1452.3
$ws.Range("I136").Value = 1369.2222
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 4107.6666
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = -1557.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 29499
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 29499
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 29499
$ws.Range("N112").Value = -32453

$ws.Range("H132").Value = 69529.22
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 69529.22
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 69529.22
$ws.Range("N132").Value = -79649.22

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 50000
$ws.Range("I23").Value = 50000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 50000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -49760
$ws.Range("N23").ClearContents()

$ws.Range("H27").Value = 50000
$ws.Range("I27").Value = 50000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 50000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -49808
$ws.Range("N27").ClearContents()

$ws.Range("H31").Value = 35718124
$ws.Range("I31").Value = 52634696
$ws.Range("J31").Value = 5355.5557
$ws.Range("K31").Value = 52634696
$ws.Range("L31").Value = 5355.5557
$ws.Range("M31").Value = -52634401
$ws.Range("N31").Value = -5945.5557

$ws.Range("H34").Value = 35718124
$ws.Range("I34").Value = 52634696
$ws.Range("J34").Value = 5355.5557
$ws.Range("K34").Value = 52634696
$ws.Range("L34").Value = 5355.5557
$ws.Range("M34").Value = -52634494
$ws.Range("N34").Value = -5759.5557

$ws.Range("H58").Value = 2181.682
$ws.Range("I58").Value = 1889
$ws.Range("J58").Value = 3498.75
$ws.Range("K58").Value = 1889
$ws.Range("L58").Value = 3498.75
$ws.Range("M58").Value = -1686
$ws.Range("N58").Value = -3904.75

$ws.Range("H122").Value = 1219.8235
$ws.Range("I122").Value = 1240.2307
$ws.Range("J122").Value = 1153.5
$ws.Range("K122").Value = 3720.6921
$ws.Range("L122").Value = 3460.5
$ws.Range("M122").Value = -1270.6921
$ws.Range("N122").Value = -8360.5

$ws.Range("H132").Value = 2358.4243
$ws.Range("I132").Value = 2213.04
$ws.Range("J132").Value = 2812.75
$ws.Range("K132").Value = 6639.12
$ws.Range("L132").Value = 8438.25
$ws.Range("M132").Value = -4109.12
$ws.Range("N132").Value = -13498.25

$ws.Range("H136").Value = 2181.682
$ws.Range("I136").Value = 1889
$ws.Range("J136").Value = 3498.75
$ws.Range("K136").Value = 5667
$ws.Range("L136").Value = 10496.25
$ws.Range("M136").Value = -3117
$ws.Range("N136").Value = -15596.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1257.8695
$ws.Range("I68").Value = 1332.6666
$ws.Range("J68").Value = 1231.4706
$ws.Range("K68").Value = 3997.9998
$ws.Range("L68").Value = 3694.4118
$ws.Range("M68").Value = -3186.9998
$ws.Range("N68").Value = -5316.4118

$ws.Range("H71").Value = 1257.8695
$ws.Range("I71").Value = 1332.6666
$ws.Range("J71").Value = 1231.4706
$ws.Range("K71").Value = 11993.9994
$ws.Range("L71").Value = 11083.2354
$ws.Range("M71").Value = -7937.999400000001
$ws.Range("N71").Value = -19195.2354

$ws.Range("H92").Value = 840
$ws.Range("I92").Value = 875
$ws.Range("J92").Value = 700
$ws.Range("K92").Value = 2625
$ws.Range("L92").Value = 2100
$ws.Range("M92").Value = -1377
$ws.Range("N92").Value = -4596

$ws.Range("H102").Value = 50000
$ws.Range("I102").Value = 50000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 150000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -147566
$ws.Range("N102").ClearContents()

$ws.Range("H107").Value = 637.53845
$ws.Range("I107").Value = 1188.25
$ws.Range("J107").Value = 392.77777
$ws.Range("K107").Value = 3564.75
$ws.Range("L107").Value = 1178.33331
$ws.Range("M107").Value = -1644.75
$ws.Range("N107").Value = -5018.33331

$ws.Range("H113").Value = 680.6957
$ws.Range("I113").Value = 552.5
$ws.Range("J113").Value = 707.6842
$ws.Range("K113").Value = 1657.5
$ws.Range("L113").Value = 2123.0526
$ws.Range("M113").Value = 512.5
$ws.Range("N113").Value = -6463.0526

$ws.Range("H118").Value = 2153.4
$ws.Range("I118").Value = 1039.8572
$ws.Range("J118").Value = 3127.75
$ws.Range("K118").Value = 3119.5716
$ws.Range("L118").Value = 9383.25
$ws.Range("M118").Value = -1876.5716
$ws.Range("N118").Value = -11869.25

$ws.Range("H131").Value = 893.9400000000001
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 893.9400000000001
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2681.82
$ws.Range("N131").Value = -12761.82
$ws.Range("M131").ClearContents()

$ws.Range("H132").Value = 1822.5883
$ws.Range("I132").Value = 1082.5
$ws.Range("J132").Value = 2226.2727
$ws.Range("K132").Value = 9742.5
$ws.Range("L132").Value = 20036.4543
$ws.Range("M132").Value = -7212.5
$ws.Range("N132").Value = -25096.4543

$ws.Range("H133").Value = 3896.875
$ws.Range("I133").Value = 1776.6666
$ws.Range("J133").Value = 5169
$ws.Range("K133").Value = 5329.9998
$ws.Range("L133").Value = 15507
$ws.Range("M133").Value = -269.9997999999996
$ws.Range("N133").Value = -25627

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 10092.526
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 10092.526
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 10092.526
$ws.Range("N123").Value = -14992.526

$ws.Range("H126").Value = 2234.739
$ws.Range("I126").Value = 2066
$ws.Range("J126").Value = 2842.2
$ws.Range("K126").Value = 6198
$ws.Range("L126").Value = 8526.599999999999
$ws.Range("M126").Value = -3728

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9572.315000000001
$ws.Range("I2").Value = 2500
$ws.Range("J2").Value = 9965.223
$ws.Range("K2").Value = 2500
$ws.Range("L2").Value = 9965.223
$ws.Range("M2").Value = -2388
$ws.Range("N2").Value = -10189.223

$ws.Range("H5").Value = 14185.714
$ws.Range("I5").Value = 14000
$ws.Range("J5").Value = 14433.333
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 14433.333
$ws.Range("M5").Value = -13887

$ws.Range("H96").Value = 29196
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 29196
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 29196
$ws.Range("N96").Value = -34688

$ws.Range("H132").Value = 4726.4863
$ws.Range("I132").Value = 4585.294
$ws.Range("J132").Value = 6326.6665
$ws.Range("K132").Value = 13755.882
$ws.Range("L132").Value = 18979.9995
$ws.Range("M132").Value = -11225.882

$ws.Range("H137").Value = 83910
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 83910
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 83910
$ws.Range("N137").Value = -94110

$ws.Range("H138").Value = 34940
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 34940
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 34940
$ws.Range("N138").Value = -45220

$ws.Range("H141").Value = 66220
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 66220
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 66220
$ws.Range("N141").Value = -76580

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 712.8570999999999
$ws.Range("I96").Value = 720
$ws.Range("J96").Value = 695
$ws.Range("K96").Value = 720
$ws.Range("L96").Value = 656.6667
$ws.Range("M96").Value = 653
$ws.Range("N96").Value = -3441

$ws.Range("H133").Value = 54001
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 54001
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 54001
$ws.Range("N133").Value = -64121
